$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.537.08"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "3.071.54"
$ws.Range("E3").Value = "  -2.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("E6").Value = "  +4.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.548"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.33%  "
$ws.Range("D9").Value = "3.076.34"
$ws.Range("E9").Value = "  -2.02%  "
$ws.Range("E10").Value = "  -3.78%  "
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.462"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("E13").Value = "  -2.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("E15").Value = "  -1.71%  "
$ws.Range("D16").Value = "3.580.88"
$ws.Range("E16").Value = "  -2.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").Value = "63.510.99"
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("D19").Value = "3.076.38"
$ws.Range("E19").Value = "  -2.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "476.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.94%  "
$ws.Range("E21").Value = "  +2.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.720"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.59%  "
$ws.Range("E23").Value = "  +0.88%  "
$ws.Range("E24").Value = "  +2.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.46%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.95%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.01%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.998"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("E33").Value = "  +3.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.82%  "
$ws.Range("E35").Value = "  +2.66%  "
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.14"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.09%  "
$ws.Range("E39").Value = "  -3.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("E41").Value = "  -1.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "448.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.27%  "
$ws.Range("E43").Value = "  -2.29%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0365"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.23%  "
$ws.Range("B45").Value = "Arweave"
$ws.Range("C45").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "40.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.44%  "
$ws.Range("E46").Value = "  +2.58%  "
$ws.Range("D47").Value = "2.812.84"
$ws.Range("E47").Value = "  -3.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.37%  "
$ws.Range("E51").Value = "  +0.56%  "
